$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 12, shifting old rows 12-20 down to 13-21 ---
$ws.Rows("12:12").Insert()

# --- Populate new row 12 (dct:creator / ORCID) ---
$ws.Cells.Item(12,1).Value = "dct:creator"
$ws.Cells.Item(12,2).Value = "0000-0003-2195-3997"

# --- Update the "dct:modified" timestamp in row 16 ---
$ws.Cells.Item(16,2).Value = "2022-06-20T11:30:28+00:00"

# --- Add newly-populated skos:broader ("memory") links for rows 19 and 20 ---
$ws.Cells.Item(19,6).Value = "memory"
$ws.Cells.Item(20,6).Value = "memory"

# --- Force text format for date-like columns L (dct:modified^^xsd:date) and M (dct:created^^xsd:date) ---
# so Excel does not auto-convert "YYYY-MM-DD" strings into date serial numbers ---
$ws.Range("L22:M33").NumberFormat = "@"

# --- Append new vocabulary term rows 22-33 ---
# Row 22: vocab:1003 / aging
$ws.Cells.Item(22,1).Value = "vocab:1003"
$ws.Cells.Item(22,2).Value = "aging"
$ws.Cells.Item(22,12).Value = "2022-06-20"
$ws.Cells.Item(22,13).Value = "2022-06-20"
$ws.Cells.Item(22,14).Value = "0000-0003-0152-1441"

# Row 23: vocab:1004 / self-reported memory
$ws.Cells.Item(23,1).Value = "vocab:1004"
$ws.Cells.Item(23,2).Value = "self-reported memory"
$ws.Cells.Item(23,6).Value = "memory"
$ws.Cells.Item(23,12).Value = "2022-06-20"
$ws.Cells.Item(23,13).Value = "2022-06-20"
$ws.Cells.Item(23,14).Value = "0000-0003-0152-1441"

# Row 24: vocab:1005 / memory
$ws.Cells.Item(24,1).Value = "vocab:1005"
$ws.Cells.Item(24,2).Value = "memory"
$ws.Cells.Item(24,12).Value = "2022-06-20"
$ws.Cells.Item(24,13).Value = "2022-06-20"
$ws.Cells.Item(24,14).Value = "0000-0003-0152-1441"

# Row 25: vocab:1006 / cognition
$ws.Cells.Item(25,1).Value = "vocab:1006"
$ws.Cells.Item(25,2).Value = "cognition"
$ws.Cells.Item(25,12).Value = "2022-06-20"
$ws.Cells.Item(25,13).Value = "2022-06-20"
$ws.Cells.Item(25,14).Value = "0000-0003-2195-3997"

# Row 26: vocab:1007 / working memory
$ws.Cells.Item(26,1).Value = "vocab:1007"
$ws.Cells.Item(26,2).Value = "working memory"
$ws.Cells.Item(26,6).Value = "memory"
$ws.Cells.Item(26,12).Value = "2022-06-20"
$ws.Cells.Item(26,13).Value = "2022-06-20"
$ws.Cells.Item(26,14).Value = "0000-0003-2195-3997"

# Row 27: vocab:1008 / attention
$ws.Cells.Item(27,1).Value = "vocab:1008"
$ws.Cells.Item(27,2).Value = "attention"
$ws.Cells.Item(27,12).Value = "2022-06-20"
$ws.Cells.Item(27,13).Value = "2022-06-20"
$ws.Cells.Item(27,14).Value = "0000-0003-2195-3997"

# Row 28: vocab:1009 / maintenance mechanism
$ws.Cells.Item(28,1).Value = "vocab:1009"
$ws.Cells.Item(28,2).Value = "maintenance mechanism"
$ws.Cells.Item(28,12).Value = "2022-06-20"
$ws.Cells.Item(28,13).Value = "2022-06-20"
$ws.Cells.Item(28,14).Value = "0000-0003-2195-3997"

# Row 29: vocab:1010 / cognitive regulation
$ws.Cells.Item(29,1).Value = "vocab:1010"
$ws.Cells.Item(29,2).Value = "cognitive regulation"
$ws.Cells.Item(29,6).Value = "regulation"
$ws.Cells.Item(29,12).Value = "2022-06-20"
$ws.Cells.Item(29,13).Value = "2022-06-20"
$ws.Cells.Item(29,14).Value = "0000-0003-2195-3997"

# Row 30: vocab:1011 / behavioral regulation
$ws.Cells.Item(30,1).Value = "vocab:1011"
$ws.Cells.Item(30,2).Value = "behavioral regulation"
$ws.Cells.Item(30,6).Value = "regulation"
$ws.Cells.Item(30,12).Value = "2022-06-20"
$ws.Cells.Item(30,13).Value = "2022-06-20"
$ws.Cells.Item(30,14).Value = "0000-0003-2195-3997"

# Row 31: vocab:1012 / emotion regulation
$ws.Cells.Item(31,1).Value = "vocab:1012"
$ws.Cells.Item(31,2).Value = "emotion regulation"
$ws.Cells.Item(31,6).Value = "regulation"
$ws.Cells.Item(31,12).Value = "2022-06-20"
$ws.Cells.Item(31,13).Value = "2022-06-20"
$ws.Cells.Item(31,14).Value = "0000-0003-2195-3997"

# Row 32: vocab:1013 / regulation
$ws.Cells.Item(32,1).Value = "vocab:1013"
$ws.Cells.Item(32,2).Value = "regulation"
$ws.Cells.Item(32,12).Value = "2022-06-20"
$ws.Cells.Item(32,13).Value = "2022-06-20"
$ws.Cells.Item(32,14).Value = "0000-0003-2195-3997"

# Row 33: vocab:1014 / sensitivity
$ws.Cells.Item(33,1).Value = "vocab:1014"
$ws.Cells.Item(33,2).Value = "sensitivity"
$ws.Cells.Item(33,12).Value = "2022-06-20"
$ws.Cells.Item(33,13).Value = "2022-06-20"
$ws.Cells.Item(33,14).Value = "0000-0003-2195-3997"

